$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Roster was re-shuffled: "Taurean Prince" (Milwaukee Bucks) was removed and
# "Aaron Wiggins" (Oklahoma City Thunder) was inserted, shifting several
# players' Position/Team pairings. Columns A (names) stayed put except for
# the following rows whose Position/Team (and in some places the player
# itself) changed to reflect the new roster order.

$ws.Range("A8").Value = "Ayo Dosunmu"
$ws.Range("B8").Value = "PG,SG,SF"
$ws.Range("C8").Value = "Chicago Bulls"

$ws.Range("A9").Value = "Pascal Siakam"
$ws.Range("B9").Value = "SF,PF,C"
$ws.Range("C9").Value = "Indiana Pacers"

$ws.Range("A11").Value = "Rudy Gobert"
$ws.Range("B11").Value = "C"
$ws.Range("C11").Value = "Minnesota Timberwolves"

$ws.Range("A12").Value = "Jalen Green"
$ws.Range("B12").Value = "PG,SG"
$ws.Range("C12").Value = "Houston Rockets"

$ws.Range("A13").Value = "Aaron Wiggins"
$ws.Range("B13").Value = "SG,SF"
$ws.Range("C13").Value = "Oklahoma City Thunder"

$ws.Range("A14").Value = "Chet Holmgren"
$ws.Range("B14").Value = "PF,C"
$ws.Range("C14").Value = "Oklahoma City Thunder"

$ws.Range("A15").Value = "Jaylen Brown"
$ws.Range("B15").Value = "SG,SF"
$ws.Range("C15").Value = "Boston Celtics"
